# Regenerate the localization handback status report:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn / de-de columns) and on each language sheet.
#  - The zh-cn / de-de "Latest Handback DateTime" values are refreshed.
#  - The (now resolved) "Error Detail" explaining a stale handback file is cleared.
#  - The Status / Error Detail columns are widened / narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# -- Overview sheet -----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# -- zh-cn sheet ----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-19 14:55:41"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# -- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-19 14:55:48"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
